$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the results of the latest "pelada" (match day) to the Jogadores log.
$names = New-Object "object[,]" 45,1
$stats = New-Object "object[,]" 45,10

$names[0,0] = "Digão"
$stats[0,0] = 4
$stats[0,1] = 3
$stats[0,2] = 3
$stats[0,3] = 1
$stats[0,4] = 1
$stats[0,5] = 1
$stats[0,6] = 0
$stats[0,7] = 0
$stats[0,8] = 0
$stats[0,9] = 0

$names[1,0] = "Marcelão"
$stats[1,0] = 4
$stats[1,1] = 3
$stats[1,2] = 3
$stats[1,3] = 2
$stats[1,4] = 1
$stats[1,5] = 1
$stats[1,6] = 0
$stats[1,7] = 0
$stats[1,8] = 0
$stats[1,9] = 0

$names[2,0] = "David"
$stats[2,0] = 4
$stats[2,1] = 3
$stats[2,2] = 3
$stats[2,3] = 0
$stats[2,4] = 1
$stats[2,5] = 1
$stats[2,6] = 0
$stats[2,7] = 0
$stats[2,8] = 0
$stats[2,9] = 0

$names[3,0] = "Juscielio"
$stats[3,0] = 4
$stats[3,1] = 3
$stats[3,2] = 3
$stats[3,3] = 2
$stats[3,4] = 1
$stats[3,5] = 1
$stats[3,6] = 0
$stats[3,7] = 0
$stats[3,8] = 0
$stats[3,9] = 0

$names[4,0] = "Eduardo"
$stats[4,0] = 4
$stats[4,1] = 3
$stats[4,2] = 3
$stats[4,3] = 5
$stats[4,4] = 1
$stats[4,5] = 1
$stats[4,6] = 0
$stats[4,7] = 0
$stats[4,8] = 0
$stats[4,9] = 0

$names[5,0] = "Fabinho"
$stats[5,0] = 4
$stats[5,1] = 3
$stats[5,2] = 3
$stats[5,3] = 0
$stats[5,4] = 1
$stats[5,5] = 0
$stats[5,6] = 0
$stats[5,7] = 0
$stats[5,8] = 0
$stats[5,9] = 0

$names[6,0] = "Caio"
$stats[6,0] = 4
$stats[6,1] = 3
$stats[6,2] = 3
$stats[6,3] = 5
$stats[6,4] = 1
$stats[6,5] = 0
$stats[6,6] = 0
$stats[6,7] = 0
$stats[6,8] = 0
$stats[6,9] = 0

$names[7,0] = "Jorge"
$stats[7,0] = 4
$stats[7,1] = 3
$stats[7,2] = 3
$stats[7,3] = 1
$stats[7,4] = 1
$stats[7,5] = 0
$stats[7,6] = 0
$stats[7,7] = 0
$stats[7,8] = 0
$stats[7,9] = 0

$names[8,0] = "Cabeleira"
$stats[8,0] = 4
$stats[8,1] = 3
$stats[8,2] = 3
$stats[8,3] = 0
$stats[8,4] = 1
$stats[8,5] = 0
$stats[8,6] = 0
$stats[8,7] = 0
$stats[8,8] = 0
$stats[8,9] = 0

$names[9,0] = "Peixe"
$stats[9,0] = 4
$stats[9,1] = 3
$stats[9,2] = 3
$stats[9,3] = 2
$stats[9,4] = 1
$stats[9,5] = 0
$stats[9,6] = 0
$stats[9,7] = 0
$stats[9,8] = 0
$stats[9,9] = 0

$names[10,0] = "Coxinha"
$stats[10,0] = 3
$stats[10,1] = 1
$stats[10,2] = 4
$stats[10,3] = 2
$stats[10,4] = 1
$stats[10,5] = 0
$stats[10,6] = 0
$stats[10,7] = 0
$stats[10,8] = 0
$stats[10,9] = 0

$names[11,0] = "Guinha"
$stats[11,0] = 3
$stats[11,1] = 1
$stats[11,2] = 4
$stats[11,3] = 0
$stats[11,4] = 1
$stats[11,5] = 0
$stats[11,6] = 0
$stats[11,7] = 0
$stats[11,8] = 0
$stats[11,9] = 0

$names[12,0] = "Nenzinho"
$stats[12,0] = 3
$stats[12,1] = 1
$stats[12,2] = 4
$stats[12,3] = 2
$stats[12,4] = 1
$stats[12,5] = 0
$stats[12,6] = 0
$stats[12,7] = 0
$stats[12,8] = 0
$stats[12,9] = 0

$names[13,0] = "Philipe"
$stats[13,0] = 3
$stats[13,1] = 1
$stats[13,2] = 4
$stats[13,3] = 1
$stats[13,4] = 1
$stats[13,5] = 0
$stats[13,6] = 0
$stats[13,7] = 0
$stats[13,8] = 0
$stats[13,9] = 0

$names[14,0] = "Vander"
$stats[14,0] = 3
$stats[14,1] = 1
$stats[14,2] = 4
$stats[14,3] = 2
$stats[14,4] = 1
$stats[14,5] = 0
$stats[14,6] = 0
$stats[14,7] = 0
$stats[14,8] = 0
$stats[14,9] = 0

$names[15,0] = "Corinthiano"
$stats[15,0] = 2
$stats[15,1] = 3
$stats[15,2] = 3
$stats[15,3] = 1
$stats[15,4] = 1
$stats[15,5] = 0
$stats[15,6] = 1
$stats[15,7] = 0
$stats[15,8] = 0
$stats[15,9] = 0

$names[16,0] = "Athos"
$stats[16,0] = 2
$stats[16,1] = 3
$stats[16,2] = 3
$stats[16,3] = 4
$stats[16,4] = 1
$stats[16,5] = 0
$stats[16,6] = 1
$stats[16,7] = 0
$stats[16,8] = 0
$stats[16,9] = 0

$names[17,0] = "Miqueias"
$stats[17,0] = 2
$stats[17,1] = 3
$stats[17,2] = 3
$stats[17,3] = 0
$stats[17,4] = 1
$stats[17,5] = 0
$stats[17,6] = 1
$stats[17,7] = 0
$stats[17,8] = 0
$stats[17,9] = 0

$names[18,0] = "Leandrinho"
$stats[18,0] = 2
$stats[18,1] = 3
$stats[18,2] = 3
$stats[18,3] = 1
$stats[18,4] = 1
$stats[18,5] = 0
$stats[18,6] = 1
$stats[18,7] = 0
$stats[18,8] = 0
$stats[18,9] = 0

$names[19,0] = "Eder"
$stats[19,0] = 2
$stats[19,1] = 3
$stats[19,2] = 3
$stats[19,3] = 0
$stats[19,4] = 1
$stats[19,5] = 0
$stats[19,6] = 1
$stats[19,7] = 0
$stats[19,8] = 0
$stats[19,9] = 0

$names[20,0] = "Matheus"
$stats[20,0] = 4
$stats[20,1] = 5
$stats[20,2] = 7
$stats[20,3] = 0
$stats[20,4] = 1
$stats[20,5] = 0
$stats[20,6] = 1
$stats[20,7] = 0
$stats[20,8] = 15
$stats[20,9] = 0

$names[21,0] = "Chelin"
$stats[21,0] = 7
$stats[21,1] = 5
$stats[21,2] = 4
$stats[21,3] = 0
$stats[21,4] = 1
$stats[21,5] = 1
$stats[21,6] = 0
$stats[21,7] = 0
$stats[21,8] = 11
$stats[21,9] = 0

$names[22,0] = "Coxinha"
$stats[22,0] = 7
$stats[22,1] = 1
$stats[22,2] = 3
$stats[22,3] = 4
$stats[22,4] = 1
$stats[22,5] = 1
$stats[22,6] = 0
$stats[22,7] = 0
$stats[22,8] = 0
$stats[22,9] = 0

$names[23,0] = "Marcelão"
$stats[23,0] = 7
$stats[23,1] = 1
$stats[23,2] = 3
$stats[23,3] = 1
$stats[23,4] = 1
$stats[23,5] = 1
$stats[23,6] = 0
$stats[23,7] = 0
$stats[23,8] = 0
$stats[23,9] = 0

$names[24,0] = "Philipe"
$stats[24,0] = 7
$stats[24,1] = 1
$stats[24,2] = 3
$stats[24,3] = 3
$stats[24,4] = 1
$stats[24,5] = 1
$stats[24,6] = 0
$stats[24,7] = 0
$stats[24,8] = 0
$stats[24,9] = 0

$names[25,0] = "Peixe"
$stats[25,0] = 7
$stats[25,1] = 1
$stats[25,2] = 3
$stats[25,3] = 1
$stats[25,4] = 1
$stats[25,5] = 1
$stats[25,6] = 0
$stats[25,7] = 0
$stats[25,8] = 0
$stats[25,9] = 0

$names[26,0] = "mateus"
$stats[26,0] = 7
$stats[26,1] = 1
$stats[26,2] = 3
$stats[26,3] = 3
$stats[26,4] = 1
$stats[26,5] = 1
$stats[26,6] = 0
$stats[26,7] = 0
$stats[26,8] = 0
$stats[26,9] = 0

$names[27,0] = "Fernando"
$stats[27,0] = 5
$stats[27,1] = 1
$stats[27,2] = 3
$stats[27,3] = 2
$stats[27,4] = 1
$stats[27,5] = 0
$stats[27,6] = 0
$stats[27,7] = 0
$stats[27,8] = 0
$stats[27,9] = 0

$names[28,0] = "Douglas"
$stats[28,0] = 5
$stats[28,1] = 1
$stats[28,2] = 3
$stats[28,3] = 2
$stats[28,4] = 1
$stats[28,5] = 0
$stats[28,6] = 0
$stats[28,7] = 0
$stats[28,8] = 0
$stats[28,9] = 0

$names[29,0] = "David"
$stats[29,0] = 5
$stats[29,1] = 1
$stats[29,2] = 3
$stats[29,3] = 0
$stats[29,4] = 1
$stats[29,5] = 0
$stats[29,6] = 0
$stats[29,7] = 0
$stats[29,8] = 0
$stats[29,9] = 0

$names[30,0] = "Leandrinho"
$stats[30,0] = 5
$stats[30,1] = 1
$stats[30,2] = 3
$stats[30,3] = 3
$stats[30,4] = 1
$stats[30,5] = 0
$stats[30,6] = 0
$stats[30,7] = 0
$stats[30,8] = 0
$stats[30,9] = 0

$names[31,0] = "Digão"
$stats[31,0] = 5
$stats[31,1] = 1
$stats[31,2] = 3
$stats[31,3] = 2
$stats[31,4] = 1
$stats[31,5] = 0
$stats[31,6] = 0
$stats[31,7] = 0
$stats[31,8] = 0
$stats[31,9] = 0

$names[32,0] = "Jorge"
$stats[32,0] = 0
$stats[32,1] = 0
$stats[32,2] = 6
$stats[32,3] = 0
$stats[32,4] = 1
$stats[32,5] = 0
$stats[32,6] = 1
$stats[32,7] = 0
$stats[32,8] = 0
$stats[32,9] = 0

$names[33,0] = "Fabinho"
$stats[33,0] = 0
$stats[33,1] = 0
$stats[33,2] = 6
$stats[33,3] = 0
$stats[33,4] = 1
$stats[33,5] = 0
$stats[33,6] = 1
$stats[33,7] = 0
$stats[33,8] = 0
$stats[33,9] = 0

$names[34,0] = "Athos"
$stats[34,0] = 0
$stats[34,1] = 0
$stats[34,2] = 6
$stats[34,3] = 0
$stats[34,4] = 1
$stats[34,5] = 0
$stats[34,6] = 1
$stats[34,7] = 0
$stats[34,8] = 0
$stats[34,9] = 0

$names[35,0] = "Eder"
$stats[35,0] = 0
$stats[35,1] = 0
$stats[35,2] = 6
$stats[35,3] = 0
$stats[35,4] = 1
$stats[35,5] = 0
$stats[35,6] = 1
$stats[35,7] = 0
$stats[35,8] = 0
$stats[35,9] = 0

$names[36,0] = "Eduardo"
$stats[36,0] = 0
$stats[36,1] = 0
$stats[36,2] = 6
$stats[36,3] = 1
$stats[36,4] = 1
$stats[36,5] = 0
$stats[36,6] = 1
$stats[36,7] = 0
$stats[36,8] = 0
$stats[36,9] = 0

$names[37,0] = "Juscielio"
$stats[37,0] = 3
$stats[37,1] = 2
$stats[37,2] = 3
$stats[37,3] = 1
$stats[37,4] = 1
$stats[37,5] = 0
$stats[37,6] = 0
$stats[37,7] = 0
$stats[37,8] = 0
$stats[37,9] = 0

$names[38,0] = "Corinthiano"
$stats[38,0] = 3
$stats[38,1] = 2
$stats[38,2] = 3
$stats[38,3] = 1
$stats[38,4] = 1
$stats[38,5] = 0
$stats[38,6] = 0
$stats[38,7] = 0
$stats[38,8] = 0
$stats[38,9] = 0

$names[39,0] = "Deiverson"
$stats[39,0] = 3
$stats[39,1] = 2
$stats[39,2] = 3
$stats[39,3] = 1
$stats[39,4] = 1
$stats[39,5] = 0
$stats[39,6] = 0
$stats[39,7] = 0
$stats[39,8] = 0
$stats[39,9] = 0

$names[40,0] = "Vander"
$stats[40,0] = 3
$stats[40,1] = 2
$stats[40,2] = 3
$stats[40,3] = 0
$stats[40,4] = 1
$stats[40,5] = 0
$stats[40,6] = 0
$stats[40,7] = 0
$stats[40,8] = 0
$stats[40,9] = 0

$names[41,0] = "Cabeleira"
$stats[41,0] = 3
$stats[41,1] = 2
$stats[41,2] = 3
$stats[41,3] = 5
$stats[41,4] = 1
$stats[41,5] = 0
$stats[41,6] = 0
$stats[41,7] = 0
$stats[41,8] = 0
$stats[41,9] = 0

$names[42,0] = "Matheus"
$stats[42,0] = 5
$stats[42,1] = 2
$stats[42,2] = 5
$stats[42,3] = 0
$stats[42,4] = 1
$stats[42,5] = 1
$stats[42,6] = 0
$stats[42,7] = 0
$stats[42,8] = 11
$stats[42,9] = 0

$names[43,0] = "Lucian"
$stats[43,0] = 5
$stats[43,1] = 1
$stats[43,2] = 5
$stats[43,3] = 0
$stats[43,4] = 1
$stats[43,5] = 0
$stats[43,6] = 0
$stats[43,7] = 0
$stats[43,8] = 8
$stats[43,9] = 0

$names[44,0] = "Breno"
$stats[44,0] = 5
$stats[44,1] = 1
$stats[44,2] = 5
$stats[44,3] = 0
$stats[44,4] = 1
$stats[44,5] = 0
$stats[44,6] = 1
$stats[44,7] = 0
$stats[44,8] = 9
$stats[44,9] = 0

$ws.Range("A443:A487").Value2 = $names
$ws.Range("C443:L487").Value2 = $stats

$ws.Application.Goto($ws.Range("K488"))
